$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados..." timestamp update (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 17:05"

# --- Row 11: India ---
$ws.Range("B11").Value = 186321
$ws.Range("C11").Value = 4494
$ws.Range("D11").Value = 88808
$ws.Range("E11").Value = 92244
$ws.Range("G11").Value = 84
$ws.Range("H11").Value = 5269

# --- Row 12: Alemania ---
$ws.Range("B12").Value = 183370
$ws.Range("C12").Value = 76
$ws.Range("E12").Value = 9568

# --- Row 29: Singapur ---
$ws.Range("D29").Value = 21699
$ws.Range("E29").Value = 13162

# --- Row 47: Argentina ---
$ws.Range("D47").Value = 5336
$ws.Range("E47").Value = 10348
$ws.Range("G47").Value = 2
$ws.Range("H47").Value = 530

# --- Row 54: Barein ---
$ws.Range("E54").Value = 4596
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 19

# --- Row 62: Moldavia ---
$ws.Range("B62").Value = 8251
$ws.Range("C62").Value = 153
$ws.Range("E62").Value = 3375
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 295

# --- Row 83: Grecia ---
$ws.Range("B83").Value = 2917
$ws.Range("C83").Value = 2
$ws.Range("E83").Value = 1368

# --- Rows 200/201: swap Belice <-> Santa Lucia (with their Casos activos / Muertes) ---
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

# --- Rows 213/214: swap Islas Virgenes Britanicas <-> Papua Nueva Guinea ---
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
